$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("126 Racecourse Road Public Housing Tower Flemington", 10, "old"),
    @("3175 The Bays Aged Care Facility Hastings", 16, "new"),
    @("3175 The Bays Aged Care Facility Hastings", 17, "old"),
    @("3600 Belvedere Age Care Noble Park Outbreak", 18, "old"),
    @("3600 Belvedere Age Care Noble Park Outbreak", 21, "new"),
    @("3612 BlueCross Glengowrie Outbreak", 22, "old"),
    @("3612 BlueCross Glengowrie Outbreak", 30, "new"),
    @("Australian Lamb Colac East", 13, "new"),
    @("Australian Lamb Colac East", 14, "old"),
    @("CS Square Caroline Springs Outbreak", 16, "old"),
    @("CS Square Caroline Springs Outbreak", 17, "new"),
    @("Child's Play Early Learning Centre Tarneit", 10, "old"),
    @("Child's Play Early Learning Centre Tarneit", 11, "new"),
    @("Embracia Aged Care Reservoir Outbreak", 22, "new"),
    @("Embracia Aged Care Reservoir Outbreak", 23, "old"),
    @("FedEx Station Melbourne Airport", 10, "old"),
    @("Hello Fresh Warehouse Ravenhall", 13, "new"),
    @("Inghams Enterprise Somerville Outbreak", 12, "old"),
    @("Inghams Enterprise Somerville Outbreak", 15, "new"),
    @("Kool Kidz Childcare Narre Warren", 10, "new"),
    @("Kool Kidz Childcare Narre Warren", 11, "old"),
    @("Lantmannen Unibake Australia Mordialloc", 25, "new"),
    @("Lantmannen Unibake Australia Mordialloc", 26, "old"),
    @("Launch Housing City Edge Crisis Accommodation South Melbourne", 12, "new"),
    @("Nido Early School Ascot Vale", 11, "new"),
    @("Nido Early School Ascot Vale", 15, "old"),
    @("Nido Early School Glenroy", 15, "old"),
    @("Northern Health Northern Hospital Epping Emergency Department Tier 1B", 42, "new"),
    @("Northern Health Northern Hospital Epping Emergency Department Tier 1B", 44, "old"),
    @("Northern Health The Northern Hospital Epping", 13, "old"),
    @("Northern Health The Northern Hospital Epping", 15, "new"),
    @("Robin Hood Inn Drouin West Outbreak", 32, "old"),
    @("Robin Hood Inn Drouin West Outbreak", 41, "new"),
    @("Social Gathering Warrnambool 28 Sep Outbreak", 13, "new"),
    @("Social Gathering Warrnambool 28 Sep Outbreak", 17, "old"),
    @("St Vincents Hospital Emergency Department Melbourne", 39, "old"),
    @("St Vincents Hospital Emergency Department Melbourne", 42, "new"),
    @("Target Distribution Centre Truganina Outbreak", 20, "new"),
    @("Target Distribution Centre Truganina Outbreak", 21, "old"),
    @("The Toolshed Bar Private Event Noojee", 10, "old"),
    @("The Toolshed Bar Private Event Noojee", 17, "new"),
    @("Visy Recycling Springvale", 29, "new"),
    @("Visy Recycling Springvale", 31, "old"),
    @("Werribee Mercy Hospital Emergency Department", 24, "new"),
    @("Werribee Mercy Hospital Emergency Department", 25, "old"),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
